# Updated BGR model - 2025-08-13 10:34
$wb = $excel.ActiveWorkbook

# --- "solar" sheet (sheet5.xml): shuffle the cost-class numbers in column P ---
$wsSolar = $wb.Worksheets.Item("solar")
$wsSolar.Range("P4").Value = 3
$wsSolar.Range("P5").Value = 4
$wsSolar.Range("P6").Value = 2

# --- "wind" sheet (sheet6.xml): shuffle the cost-class numbers in column P ---
$wsWind = $wb.Worksheets.Item("wind")
$wsWind.Range("P15").Value = 3
$wsWind.Range("P16").Value = 1
$wsWind.Range("P17").Value = 2

$wsWind.Range("P18").Value = 3
$wsWind.Range("P19").Value = 1
$wsWind.Range("P20").Value = 2

$wsWind.Range("P27").Value = 2
$wsWind.Range("P28").Value = 1

$wsWind.Range("P47").Value = 2
$wsWind.Range("P48").Value = 1
